$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to text format so numeric-looking strings
# (e.g. "544.24", "66.946.55") are stored as literal text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.946.55'
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').Value = '2.368.03'
$ws.Range('E3').Value = '  -4.06%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '544.24'
$ws.Range('E5').Value = '  -2.64%  '
$ws.Range('D6').Value = '155.36'
$ws.Range('E6').Value = '  -4.81%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.498'
$ws.Range('E8').Value = '  -1.10%  '
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('D11').Value = '0.321'
$ws.Range('E11').Value = '  -4.12%  '
$ws.Range('D12').Value = '4.68'
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('D13').Value = '66.872.59'
$ws.Range('E13').Value = '  -2.92%  '
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').Value = '22.41'
$ws.Range('E15').Value = '  -5.18%  '
$ws.Range('D16').Value = '10.12'
$ws.Range('E16').Value = '  -6.17%  '
$ws.Range('D17').Value = '323.15'
$ws.Range('E17').Value = '  -5.77%  '
$ws.Range('D18').Value = '6.64'
$ws.Range('E18').Value = '  -6.18%  '
$ws.Range('D19').Value = '3.69'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').Value = '1.80'
$ws.Range('E21').Value = '  -6.21%  '
$ws.Range('D22').Value = '64.98'
$ws.Range('E22').Value = '  -3.14%  '
$ws.Range('D23').Value = '3.52'
$ws.Range('E23').Value = '  -4.68%  '
$ws.Range('D24').Value = '7.77'
$ws.Range('E24').Value = '  -5.39%  '
$ws.Range('D25').Value = '0.0₃0777'
$ws.Range('E25').Value = '  -5.08%  '
$ws.Range('D26').Value = '6.88'
$ws.Range('E26').Value = '  -4.61%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '405.49'
$ws.Range('E28').Value = '  -7.76%  '
$ws.Range('D29').Value = '1.09'
$ws.Range('E29').Value = '  -4.51%  '
$ws.Range('E30').Value = '  -3.66%  '
$ws.Range('D31').Value = '158.68'
$ws.Range('E31').Value = '  +1.69%  '
$ws.Range('D32').Value = '18.92'
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').Value = '17.38'
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('E35').Value = '  -6.65%  '
$ws.Range('D36').Value = '0.286'
$ws.Range('E36').Value = '  -5.29%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '4.13'
$ws.Range('E37').Value = '  -7.49%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.42'
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('D39').Value = '1.03'
$ws.Range('E39').Value = '  -6.30%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.21'
$ws.Range('E40').Value = '  -4.43%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '125.98'
$ws.Range('E41').Value = '  -5.47%  '
$ws.Range('D42').Value = '1.90'
$ws.Range('E42').Value = '  -8.46%  '
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D44').Value = '0.461'
$ws.Range('E44').Value = '  -4.61%  '
$ws.Range('D45').Value = '0.542'
$ws.Range('E45').Value = '  -3.61%  '
$ws.Range('D46').Value = '0.0903'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').Value = '1.30'
$ws.Range('E48').Value = '  -9.65%  '
$ws.Range('D49').Value = '16.07'
$ws.Range('E49').Value = '  -5.08%  '
$ws.Range('D50').Value = '0.0419'
$ws.Range('E50').Value = '  -2.73%  '
$ws.Range('D51').Value = '0.0₆0195'
$ws.Range('E51').Value = '  -7.46%  '

# Restore column D to the default (unstyled) cell format so no stray
# style index is left referenced on these cells.
$ws.Range("D2:D51").Style = "Normal"

